# Update "想去人数" (F column) values on the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets,
# per commit "Update gh-pages to output generated at 7921097".

$wb = $excel.ActiveWorkbook

# Sheet index 1 = 展览
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F5").Value = 836
$wsExhibit.Range("F6").Value = 510
$wsExhibit.Range("F12").Value = 1996
$wsExhibit.Range("F14").Value = 995
$wsExhibit.Range("F15").Value = 2634
$wsExhibit.Range("F19").Value = 145
$wsExhibit.Range("F21").Value = 210
$wsExhibit.Range("F27").Value = 1069
$wsExhibit.Range("F29").Value = 2550
$wsExhibit.Range("F35").Value = 201
$wsExhibit.Range("F37").Value = 183

# Sheet index 4 = 全部类型
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F7").Value = 836
$wsAll.Range("F8").Value = 510
$wsAll.Range("F14").Value = 1996
$wsAll.Range("F16").Value = 995
$wsAll.Range("F18").Value = 2634
$wsAll.Range("F22").Value = 145
$wsAll.Range("F24").Value = 210
$wsAll.Range("F31").Value = 1069
$wsAll.Range("F33").Value = 2550
$wsAll.Range("F39").Value = 201
$wsAll.Range("F41").Value = 183
